$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (row 19), which drops out of the table entirely
$ws.Rows(19).Delete()

$ws.Range("A2").Value = 39765
$ws.Range("B2").Value = 2008
$ws.Range("D2").Value = 2009

$ws.Range("A3").Value = 40130
$ws.Range("B3").Value = 2009
$ws.Range("C3").Value = 0.1715429114845124
$ws.Range("D3").Value = 2010
$ws.Range("E3").ClearContents()

$ws.Range("A4").Value = 40494
$ws.Range("B4").Value = 2010
$ws.Range("C4").Value = 0.8004663283405655
$ws.Range("D4").Value = 2011
$ws.Range("E4").ClearContents()

$ws.Range("A5").Value = 40862
$ws.Range("B5").Value = 2011
$ws.Range("C5").Value = 5.253783907501819
$ws.Range("D5").Value = 2012
$ws.Range("E5").ClearContents()

$ws.Range("A6").Value = 41228
$ws.Range("B6").Value = 2012
$ws.Range("C6").Value = 3.522405026196918
$ws.Range("D6").Value = 2013
$ws.Range("E6").Value = 0.5784444854042281

$ws.Range("A7").Value = 41592
$ws.Range("B7").Value = 2013
$ws.Range("C7").Value = 1.656063945467268
$ws.Range("D7").Value = 2014
$ws.Range("E7").Value = 2.529895848567842

$ws.Range("A8").Value = 41957
$ws.Range("B8").Value = 2014
$ws.Range("C8").Value = 4.06235252733802
$ws.Range("D8").Value = 2015
$ws.Range("E8").Value = 4.060884847379076

$ws.Range("A9").Value = 42321
$ws.Range("B9").Value = 2015
$ws.Range("C9").Value = 3.05427116350534
$ws.Range("D9").Value = 2016
$ws.Range("E9").Value = 2.270469368501771

$ws.Range("A10").Value = 42689
$ws.Range("B10").Value = 2016
$ws.Range("C10").Value = 2.305809238174006
$ws.Range("D10").Value = 2017
$ws.Range("E10").Value = 2.467161166346266

$ws.Range("A11").Value = 43053
$ws.Range("B11").Value = 2017
$ws.Range("C11").Value = 2.509111342826809
$ws.Range("D11").Value = 2018
$ws.Range("E11").Value = 2.480855794925163

$ws.Range("A12").Value = 43418
$ws.Range("B12").Value = 2018
$ws.Range("C12").Value = 3.296731496509198
$ws.Range("D12").Value = 2019
$ws.Range("E12").Value = 3.221757900820066

$ws.Range("A13").Value = 43783
$ws.Range("B13").Value = 2019
$ws.Range("C13").Value = 2.861315725866587
$ws.Range("D13").Value = 2020
$ws.Range("E13").Value = 2.631992339577627

$ws.Range("A14").Value = 44159
$ws.Range("B14").Value = 2020
$ws.Range("C14").Value = 1.790319754067715
$ws.Range("D14").Value = 2021
$ws.Range("E14").Value = 2.153309886824961

$ws.Range("A15").Value = 44525
$ws.Range("B15").Value = 2021
$ws.Range("C15").Value = 2.339531676162721
$ws.Range("D15").Value = 2022
$ws.Range("E15").Value = 4.667362054855917

$ws.Range("A16").Value = 44890
$ws.Range("B16").Value = 2022
$ws.Range("C16").Value = 4.834496776263886
$ws.Range("D16").Value = 2023
$ws.Range("E16").Value = 3.305715257492858

$ws.Range("A17").Value = 45254
$ws.Range("B17").Value = 2023
$ws.Range("C17").Value = 2.798216547494237
$ws.Range("D17").Value = 2024
$ws.Range("E17").Value = 1.757655717321982

$ws.Range("A18").Value = 45618
$ws.Range("B18").Value = 2024
$ws.Range("C18").Value = 1.530879676868468
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = 2.159361127638926
